# Update "想去人数" (F column) values on the "展览" and "全部类型" sheets
# to reflect refreshed counts from the gh-pages data regeneration.

$wb = $excel.ActiveWorkbook

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)

    if ($sheetName -eq "展览") {
        $ws.Range("F3").Value = 5524
        $ws.Range("F4").Value = 44
        $ws.Range("F6").Value = 404
        $ws.Range("F9").Value = 4366
        $ws.Range("F10").Value = 783
        $ws.Range("F11").Value = 813
        $ws.Range("F14").Value = 125
        $ws.Range("F17").Value = 16
        $ws.Range("F18").Value = 121
        $ws.Range("F19").Value = 605
        $ws.Range("F21").Value = 180
        $ws.Range("F22").Value = 1128
        $ws.Range("F23").Value = 15
        $ws.Range("F24").Value = 2757
        $ws.Range("F25").Value = 439
        $ws.Range("F26").Value = 288
    }
    elseif ($sheetName -eq "全部类型") {
        $ws.Range("F3").Value = 5524
        $ws.Range("F4").Value = 44
        $ws.Range("F6").Value = 404
        $ws.Range("F9").Value = 4366
        $ws.Range("F10").Value = 783
        $ws.Range("F11").Value = 813
        $ws.Range("F14").Value = 125
        $ws.Range("F17").Value = 16
        $ws.Range("F18").Value = 121
        $ws.Range("F19").Value = 605
        $ws.Range("F22").Value = 180
        $ws.Range("F23").Value = 1128
        $ws.Range("F24").Value = 15
        $ws.Range("F25").Value = 2757
        $ws.Range("F26").Value = 439
        $ws.Range("F27").Value = 288
    }
}

$wb.Save()
